$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as row 331 (Fruta / Piña, Agro
# Chillán). Every existing row from 331 down shifts to row+1, so insert a
# fresh row at position 331 first, then populate it with the new record.
$ws.Rows.Item(331).Insert()

$ws.Cells.Item(331, 1).Value = 7
$ws.Cells.Item(331, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(331, 3).Value = "Ñuble"
$ws.Cells.Item(331, 4).Value = 45194
$ws.Cells.Item(331, 4).NumberFormat = $ws.Cells.Item(332, 4).NumberFormat
$ws.Cells.Item(331, 5).Value = 16
$ws.Cells.Item(331, 6).Value = "Fruta"
$ws.Cells.Item(331, 7).Value = 100108
$ws.Cells.Item(331, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(331, 9).Value = 100108005
$ws.Cells.Item(331, 10).Value = "Piña"
$ws.Cells.Item(331, 11).Value = "Caramelo"
$ws.Cells.Item(331, 12).Value = "Segunda"
$ws.Cells.Item(331, 13).Value = 80
$ws.Cells.Item(331, 14).Value = 21000
$ws.Cells.Item(331, 15).Value = 21000
$ws.Cells.Item(331, 16).Value = 21000
$ws.Cells.Item(331, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(331, 18).Value = "Ecuador"
$ws.Cells.Item(331, 19).Value = 1500
$ws.Cells.Item(331, 20).Value = 14
